$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new BTC-change-influenced data point
$ws.Range("A3").Value = 42941
$ws.Range("A3").NumberFormat = "mm-dd-yy"

$ws.Range("B3").Value = 0.00018118

# Update selection to reflect the new active cell
[void]$ws.Range("C3").Select()
